$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.Value = "'66.045.43"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -1.14%  '

$c = $ws.Cells.Item(3, 4)
$c.Value = "'3.762.07"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +1.77%  '

$ws.Cells.Item(4, 5).Value = '  -0.17%  '

$c = $ws.Cells.Item(5, 4)
$c.Value = "'408.33"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -2.78%  '

$c = $ws.Cells.Item(6, 4)
$c.Value = "'132.23"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +1.62%  '

$c = $ws.Cells.Item(7, 4)
$c.Value = "'3.752.32"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +1.71%  '

$ws.Cells.Item(8, 5).Value = '  -5.24%  '

$ws.Cells.Item(9, 5).Value = '  -0.01%  '

$c = $ws.Cells.Item(10, 4)
$c.Value = "'0.728"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -4.31%  '

$c = $ws.Cells.Item(11, 4)
$c.Value = "'0.165"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -8.87%  '

$c = $ws.Cells.Item(12, 4)
$c.Value = "'0.0000356"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -7.93%  '

$c = $ws.Cells.Item(13, 4)
$c.Value = "'41.08"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -3.76%  '

$c = $ws.Cells.Item(14, 4)
$c.Value = "'4.353.47"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +1.39%  '

$c = $ws.Cells.Item(15, 4)
$c.Value = "'9.95"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -4.29%  '

$c = $ws.Cells.Item(16, 4)
$c.Value = "'14.81"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +13.82%  '

$ws.Cells.Item(17, 5).Value = '  -1.35%  '

$c = $ws.Cells.Item(18, 4)
$c.Value = "'3.761.59"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +1.03%  '

$c = $ws.Cells.Item(19, 4)
$c.Value = "'19.40"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -5.38%  '

$c = $ws.Cells.Item(20, 4)
$c.Value = "'66.250.09"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -0.96%  '

$ws.Cells.Item(21, 5).Value = '  -5.15%  '

$c = $ws.Cells.Item(22, 4)
$c.Value = "'411.57"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -8.14%  '

$c = $ws.Cells.Item(23, 4)
$c.Value = "'14.31"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -10.06%  '

$c = $ws.Cells.Item(24, 4)
$c.Value = "'85.05"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -4.90%  '

$ws.Cells.Item(25, 5).Value = '  -2.31%  '

$c = $ws.Cells.Item(26, 4)
$c.Value = "'5.71"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +14.20%  '

$c = $ws.Cells.Item(27, 4)
$c.Value = "'35.81"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -4.77%  '

$ws.Cells.Item(28, 5).Value = '  -5.86%  '

$c = $ws.Cells.Item(29, 4)
$c.Value = "'9.34"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -8.96%  '

$c = $ws.Cells.Item(30, 4)
$c.Value = "'731.78"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +11.78%  '

$ws.Cells.Item(31, 2).Value = 'Hedera'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(31, 4)
$c.Value = "'0.121"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -0.57%  '

$ws.Cells.Item(32, 2).Value = 'Cosmos'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$c = $ws.Cells.Item(32, 4)
$c.Value = "'12.32"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -1.73%  '

$ws.Cells.Item(33, 2).Value = 'Toncoin'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$c = $ws.Cells.Item(33, 4)
$c.Value = "'2.73"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +1.56%  '

$ws.Cells.Item(34, 2).Value = 'RenderToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Cells.Item(34, 4)
$c.Value = "'7.39"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +1.37%  '

$ws.Cells.Item(35, 5).Value = '  -6.19%  '

$c = $ws.Cells.Item(36, 4)
$c.Value = "'38.92"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -4.74%  '

$ws.Cells.Item(37, 5).Value = '  +0.08%  '

$c = $ws.Cells.Item(38, 4)
$c.Value = "'55.04"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -3.65%  '

$c = $ws.Cells.Item(39, 4)
$c.Value = "'0.0₃0736"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -0.74%  '

$c = $ws.Cells.Item(40, 4)
$c.Value = "'0.0459"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -6.29%  '

$ws.Cells.Item(41, 5).Value = '  -14.31%  '

$c = $ws.Cells.Item(42, 4)
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -0.02%  '

$ws.Cells.Item(43, 5).Value = '  -8.80%  '

$c = $ws.Cells.Item(44, 4)
$c.Value = "'27.01"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -5.95%  '

$ws.Cells.Item(45, 2).Value = 'Monero'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$c = $ws.Cells.Item(45, 4)
$c.Value = "'146.00"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -1.37%  '

$ws.Cells.Item(46, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Cells.Item(46, 4)
$c.Value = "'3.14"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +18.10%  '

$ws.Cells.Item(47, 5).Value = '  -3.72%  '

$c = $ws.Cells.Item(48, 4)
$c.Value = "'2.06"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -1.32%  '

$ws.Cells.Item(49, 5).Value = '  -0.85%  '

$ws.Cells.Item(50, 5).Value = '  -2.39%  '

$c = $ws.Cells.Item(51, 4)
$c.Value = "'2.79"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -3.83%  '
